$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the 148 (merge-sort / recursive) problem row as done by filling
# in the previously empty E3 cell.
$ws.Range("E3").Value = "done"

# The sheet's recorded selection moves from the old C13 to E4.
$ws.Range("E4").Select()
